# feat: add 2022-Q1 data
#
# Before:  Sheet1 "2021-Q4" (fund holdings), Sheet2 "总计" (totals)
# After:   Sheet1 "2021-Q4" (unchanged), Sheet2 "2022-Q1" (new fund holdings,
#          reusing the old "总计" sheet/rId), Sheet3 "总计" (new totals table
#          appended at the end, with a 2022-Q1 row added).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "2021-Q4" - used as a style template
$ws2 = $wb.Worksheets.Item(2)   # currently "总计" -> becomes "2022-Q1"

# ---------------------------------------------------------------------
# 1) Turn the old "总计" sheet into the new "2022-Q1" holdings sheet
# ---------------------------------------------------------------------
$ws2.Name = "2022-Q1"
$ws2.Cells.Clear()

# Copy the header-row and index-column formatting from the "2021-Q4" sheet
# so the new sheet keeps the same bold/bordered look (same style index).
$ws1.Range("B1:H1").Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)
$ws1.Range("A2:A10").Copy()
$ws2.Range("A2:A10").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws2.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

$rows = @(
    @("012284", "光大保德信健康优加混合型证券投资基金", "14.49", "91.97", "3.11", "0.4506", 9),
    @("004818", "国寿安保目标策略灵活配置混合A",       "4.06",  "36.45", "3.09", "0.1255", 2),
    @("001672", "国寿安保智慧生活股票",                 "3.56",  "85.91", "3.30", "0.1175", 3),
    @("003243", "上投摩根中国世纪灵活配置混合人民币份额（QDII）", "1.36", "84.74", "3.54", "0.0481", 3),
    @("003244", "上投摩根中国世纪灵活配置混合美元现钞（QDII）", "1.36", "84.74", "3.54", "0.0481", 3),
    @("003245", "上投摩根中国世纪灵活配置混合美元现汇（QDII）", "1.36", "84.74", "3.54", "0.0481", 3),
    @("010434", "红土创新医疗保健股票",                 "0.75",  "92.96", "4.72", "0.0354", 4),
    @("004819", "国寿安保目标策略灵活配置混合C",       "0.57",  "36.45", "3.09", "0.0176", 2),
    @("006890", "上投摩根领先优选混合",                 "0.36",  "79.50", "3.66", "0.0132", 2)
)

$r = 2
foreach ($row in $rows) {
    $ws2.Cells.Item($r, 1).Value = $r - 2
    $ws2.Cells.Item($r, 2).Value = "'" + $row[0]
    $ws2.Cells.Item($r, 3).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = "'" + $row[2]
    $ws2.Cells.Item($r, 5).Value = "'" + $row[3]
    $ws2.Cells.Item($r, 6).Value = "'" + $row[4]
    $ws2.Cells.Item($r, 7).Value = "'" + $row[5]
    $ws2.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Append a brand-new "总计" sheet at the end with the updated totals
#    (duplicate an existing sheet so sheet-level properties such as
#    <sheetPr><outlinePr .../></sheetPr> carry over, then wipe its data)
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws1.Copy($null, $last)
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "总计"
$ws3.Cells.Clear()

$ws1.Range("B1:D1").Copy()
$ws3.Range("B1:D1").PasteSpecial(-4122)
$ws1.Range("A2:A3").Copy()
$ws3.Range("A2:A3").PasteSpecial(-4122)

$ws3.Cells.Item(1, 2).Value = "日期"
$ws3.Cells.Item(1, 3).Value = "持有数量(只)"
$ws3.Cells.Item(1, 4).Value = "持有市值(亿元)"

$ws3.Cells.Item(2, 1).Value = 0
$ws3.Cells.Item(2, 2).Value = "2022-Q1"
$ws3.Cells.Item(2, 3).Value = 9
$ws3.Cells.Item(2, 4).Value = 0.9

$ws3.Cells.Item(3, 1).Value = 1
$ws3.Cells.Item(3, 2).Value = "2021-Q4"
$ws3.Cells.Item(3, 3).Value = 9
$ws3.Cells.Item(3, 4).Value = 0.71

Write-Host "Sheets:" $wb.Worksheets.Count
